# Add the "2303" minutes sheet (new meeting added after "1803"), carrying the
# same layout/styling as the previous "1803" sheet, then updating the
# date/text content for the new meeting entry.

$wb = $excel.ActiveWorkbook

# --- Locate the most recent existing minutes sheet ("1803") -----------------
$prev = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- Duplicate it (preserves column widths, styles, row heights, etc.) ------
$prev.Copy($null, $prev)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "2303"

# The source sheet had an extra "Complete stakeholder analysis" action-review
# row (row 10) that doesn't apply to this meeting - remove it, shifting
# everything below up by one row.
$ws.Rows.Item(10).Delete()

# --- Update the header date --------------------------------------------------
$ws.Range("B1").Value = "03/23/2021"

# --- Action Review: mark the two in-progress items as finished -------------
$ws.Range("C6").Value = "Finished"
$ws.Range("C7").Value = "Finished"

# --- New Actions and Owners row ---------------------------------------------
$ws.Range("A13").Value = "Everyone to complete their relevant workshop"
$ws.Range("B13").Value = "all"
$ws.Range("C13").Value = "next week"

# --- Decisions ----------------------------------------------------------------
$ws.Range("A16").Value = "Talked about timeline and what we need to be doing"
$ws.Range("A17").Value = "Start report writing now"
$ws.Range("A18").Value = "Technical development needs to be caught up"

# --- Row heights reflecting the new (re-wrapped) text -----------------------
$ws.Rows.Item(7).RowHeight = 90
$ws.Rows.Item(8).RowHeight = 180
$ws.Rows.Item(9).RowHeight = 90
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 45
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(18).RowHeight = 30

# --- View state: previous sheet scrolls/selects its full range, and is no
# longer the selected tab; new sheet becomes the active/selected tab. --------
$prev.Range("A1:C23").Select()

$ws.Activate()
$ws.Range("D16").Select()
